$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A new weekly price-report row was inserted for this market/product.
# It lands at row 390, pushing every existing row from 390 downward by one
# (old row 390 -> new row 391, ..., old row 469 -> new row 470).
$ws.Rows("390:390").Insert()

# Populate the newly inserted row 390 with the new observation.
$ws.Range("A390").Value = 4
$ws.Range("B390").Value = "Feria Lagunitas de Puerto Montt"
$ws.Range("C390").Value = "Los Lagos"
$ws.Range("D390").Value = 45244
$ws.Range("E390").Value = 10
$ws.Range("F390").Value = "Fruta"
$ws.Range("G390").Value = 100108
$ws.Range("H390").Value = "Tropicales y subtropicales"
$ws.Range("I390").Value = 100108002
$ws.Range("J390").Value = "Mango"
$ws.Range("K390").Value = "Sin especificar"
$ws.Range("L390").Value = "Primera"
$ws.Range("M390").Value = 108
$ws.Range("N390").Value = 13000
$ws.Range("O390").Value = 13000
$ws.Range("P390").Value = 13000
$ws.Range("Q390").Value = "`$/bandeja 4 kilos"
$ws.Range("R390").Value = "Brasil"
$ws.Range("S390").Value = 3250
$ws.Range("T390").Value = 4
